# Apply the "novas redes / comeco do qqplot" edit to metrics_22_7 sheet.
# Every row's metric columns (B:Q) are overwritten with a new constant set of
# values (the new ensemble-averaged metrics), and column A (the model label)
# is re-permuted across rows 2-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-label order for rows 2..26 (A column)
$newLabels = @(
    "model_22_7_0",
    "model_22_7_22",
    "model_22_7_21",
    "model_22_7_20",
    "model_22_7_19",
    "model_22_7_18",
    "model_22_7_17",
    "model_22_7_16",
    "model_22_7_15",
    "model_22_7_14",
    "model_22_7_13",
    "model_22_7_23",
    "model_22_7_12",
    "model_22_7_10",
    "model_22_7_9",
    "model_22_7_8",
    "model_22_7_7",
    "model_22_7_6",
    "model_22_7_5",
    "model_22_7_4",
    "model_22_7_3",
    "model_22_7_2",
    "model_22_7_1",
    "model_22_7_11",
    "model_22_7_24"
)

# New constant metric values shared by every data row, columns B..Q
$newMetrics = @(
    0.9649258184664765,
    0.7153214174682279,
    0.9825069475759221,
    0.9778771248369326,
    0.9823676660927012,
    0.2345412134164277,
    1.903646992214302,
    0.1970878225799122,
    0.09166454694801439,
    0.1443761847639633,
    0.4038156865157607,
    0.4842945523299097,
    1.014768076435168,
    0.5049119740505222,
    164.9002479127228,
    263.629189727047
)

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newLabels[$i]
    for ($j = 0; $j -lt $newMetrics.Length; $j++) {
        $col = $j + 2
        $ws.Cells.Item($row, $col).Value = $newMetrics[$j]
    }
}
